$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.295.95'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.60%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.313.41'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.19%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '565.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.89%  '
$ws.Range('E6').Value = '  -3.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.310.27'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.480'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.35'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.00%  '
$ws.Range('E11').Value = '  -4.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.377'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.879.39'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.312.86'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.27%  '
$ws.Range('E16').Value = '  -5.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.88'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.390.26'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.68'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -9.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '354.04'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -8.18%  '
$ws.Range('E23').Value = '  -3.83%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.444.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '69.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.35%  '
$ws.Range('E27').Value = '  -5.70%  '
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.18'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.44'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.87'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.12%  '
$ws.Range('E32').Value = '  -5.94%  '
$ws.Range('E33').Value = '  -2.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.343.56'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.63'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.10%  '
$ws.Range('E38').Value = '  -0.66%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '160.19'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.63%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.48'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.42%  '
$ws.Range('E41').Value = '  -2.45%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.38'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.743'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.83%  '
$ws.Range('E46').Value = '  -5.69%  '
$ws.Range('E47').Value = '  -5.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.32'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -8.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.74'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.868'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.25'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.83%  '

Write-Output "Updated cryptos list"